$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sun Oct 13 00:06:50 EDT 2024"
$ws.Range("B3").Value = "Sun Oct 13 00:07:05 EDT 2024"
$ws.Range("B4").Value = "Sun Oct 13 00:07:19 EDT 2024"
